$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("C17").Value = 49
$ws.Range("G17").Value = 8

# Row 18 (C18 holds a numeric-looking value stored as TEXT - use the quote
# prefix so it stays text instead of becoming a real number)
$ws.Range("C18").Value = "'14"
$ws.Range("G18").Value = 17

# Row 19 (D19 is also numeric-looking TEXT)
$ws.Range("D19").Value = "'14"
$ws.Range("G19").Value = 24

# Row 20
$ws.Range("D20").Value = "81.6% der Karten"

# Row 21
$ws.Range("A21").Value = 3
$ws.Range("B21").Value = 11

# Rows 27-31 (task cards)
$ws.Range("B27").Value = "Aufgabenteilung Lenni Laura 💩"
$ws.Range("C27").Value = 12

$ws.Range("B28").Value = "to do Laura 💩"
$ws.Range("C28").Value = 12

$ws.Range("B29").Value = "Input google site:Ausreisehandbuch 💩"
$ws.Range("C29").Value = 12

$ws.Range("B30").Value = "wiki how: Kochstellen manual"
$ws.Range("C30").Value = 12

$ws.Range("B31").Value = "Austausch mit Finanzteam über kashana"
$ws.Range("C31").Value = 12

# Row 34
$ws.Range("F34").Value = 17
$ws.Range("G34").Value = "(34.7%)"

# Rows 36-39 (aktivste Mitglieder ranking)
$ws.Range("B36").Value = "Lennard Mörsen"
$ws.Range("C36").Value = 12

$ws.Range("B37").Value = "Chris Kneip"
$ws.Range("C37").Value = 5

$ws.Range("B38").Value = "Till Schöpe"

$ws.Range("B39").Value = "Micha Landoll"
$ws.Range("C39").Value = 3

# Rows 40-42 (second ranking list)
$ws.Range("E40").Value = "Lennard Mörsen"
$ws.Range("F40").Value = 9

$ws.Range("E41").Value = "Chris Kneip"
$ws.Range("F41").Value = 4

$ws.Range("E42").Value = "Micha Landoll"
$ws.Range("F42").Value = 3

# Row 50
$ws.Range("F50").Value = 3

# Row 51
$ws.Range("B51").Value = "Lennard Mörsen"

# Row 52
$ws.Range("B52").Value = "Marie-Sophie Braun"
$ws.Range("E52").Value = "Laura Coordt"

# Row 53
$ws.Range("E53").Value = "Chris Kneip"

# Row 54
$ws.Range("E54").Value = "Lara Eisler"

# Row 60
$ws.Range("G60").Value = 24

# Row 70
$ws.Range("B70").Value = 49

# Row 71
$ws.Range("B71").Value = 7

# Row 73
$ws.Range("F73").Value = "Lennard Mörsen"

# Row 74
$ws.Range("F74").Value = "Marie-Sophie Braun"
